# Insert a new weekly observation row for "Orégano" (Vega Central Mapocho de
# Santiago) as the new row 11, pushing the existing rows 11-33 down to 12-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11-33 down one row, then add the new record in row 11.
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 44482
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = 100112029
$ws.Cells.Item(11, 7).Value = "Orégano"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 16
$ws.Cells.Item(11, 11).Value = 9000
$ws.Cells.Item(11, 12).Value = 10000
$ws.Cells.Item(11, 13).Value = 9500
$ws.Cells.Item(11, 14).Value = "`$/docena de atados"
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 3167
$ws.Cells.Item(11, 17).Value = 3
$ws.Cells.Item(11, 18).Value = "Hortaliza"
